$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15, pushing existing rows 15-21 down to 16-22.
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the new price record.
$ws.Cells.Item(15, 1).Value = 3
$ws.Cells.Item(15, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(15, 3).Value = "Coquimbo"
$ws.Cells.Item(15, 4).Value = 44452
$ws.Cells.Item(15, 5).Value = 5
$ws.Cells.Item(15, 6).Value = 100112022
$ws.Cells.Item(15, 7).Value = "Arveja Verde"
$ws.Cells.Item(15, 8).Value = "Perfection"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 70
$ws.Cells.Item(15, 11).Value = 31000
$ws.Cells.Item(15, 12).Value = 32000
$ws.Cells.Item(15, 13).Value = 31500
$ws.Cells.Item(15, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(15, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(15, 16).Value = 1260
$ws.Cells.Item(15, 17).Value = 25
$ws.Cells.Item(15, 18).Value = "Hortaliza"
